$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

# Updated values for existing rows 2-10 (columns B:G)
$data = @{
    2  = @(0.001068801475677701, 0.5121987664681366, 0.5908018752482725, 0.7686363738779688, 0.776283934891037, 51)
    3  = @(0.1076112683830704, 0.5501073581205883, 0.7154227106871884, 0.8458266434011099, 0.8474707421379045, 50)
    4  = @(0.01447167482266423, 0.6025716605263358, 0.7372352437580105, 0.8586240409853492, 0.8673987081222851, 49)
    5  = @(0.107527968228121, 0.6082289609806715, 0.7759918099762468, 0.8809039731867753, 0.883568901926184, 48)
    6  = @(0.01941581067445748, 0.5766303201776867, 0.6540589049982858, 0.8087390833873962, 0.8172468464677016, 47)
    7  = @(0.1112754059595625, 0.564041488129664, 0.7222269590935267, 0.8498393725249065, 0.8518327715854772, 46)
    8  = @(0.01597012535789634, 0.5344976713245601, 0.6191284896479551, 0.7868471831607171, 0.7955744781642835, 45)
    9  = @(0.05572434282723118, 0.6058820781989243, 0.7089638270491668, 0.8419998972975987, 0.849867014019446, 44)
    10 = @(0.06564388546255764, 0.6068818695593868, 0.7068678503580751, 0.8407543341298188, 0.8481074928832305, 43)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
    $ws.Cells.Item($r, 7).Value = $vals[5]
}

# New row 11: label "Q9" using same style as other label cells (copy A10 formatting), plus data values
$ws.Cells.Item(11, 1).Value = "Q9"
$ws.Cells.Item(10, 1).Copy()
$ws.Cells.Item(11, 1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(11, 2).Value = 0.03898798267175772
$ws.Cells.Item(11, 3).Value = 0.5674245104323993
$ws.Cells.Item(11, 4).Value = 0.5443824653879992
$ws.Cells.Item(11, 5).Value = 0.7378227872517893
$ws.Cells.Item(11, 6).Value = 0.7457231047781384
$ws.Cells.Item(11, 7).Value = 42

$wb.Save()
